# data driven test case
# - rename/repopulate the login cell on the original "TC1" sheet (username/bhanu)
# - add a new "ValidLogin" worksheet after TC1 holding a small username/password
#   data table used to drive the test case

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# TC1: keep "username" in A1, change the credential value in A2
$ws1.Range("A2").Value = "bhanu"
[void]$ws1.Range("A1:A2").Select()

# Insert the new data sheet right after TC1 and rename it
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "ValidLogin"

# Populate the data table: username/password headers + admin/pointofsale row
$ws2.Range("A1").Value = "username"
$ws2.Range("B1").Value = "password"
$ws2.Range("A2").Value = "admin"
$ws2.Range("B2").Value = "pointofsale"

# Match the view state captured in the authored workbook
$ws2.Application.ActiveWindow.Zoom = 160
[void]$ws2.Range("B3").Select()
